# Fruta / hortaliza, semanal
# Insert a new weekly record before the current row 9 (pushing existing
# rows 9-24 down to 10-25) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(9).Insert()

$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44477
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 300000000
$ws.Cells.Item(9, 7).Value = "Espárragos"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 1500
$ws.Cells.Item(9, 12).Value = 1600
$ws.Cells.Item(9, 13).Value = 1550
$ws.Cells.Item(9, 14).Value = "`$/kilo"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 1550
$ws.Cells.Item(9, 17).Value = 1
$ws.Cells.Item(9, 18).Value = "Hortaliza"
